$d = $word.ActiveDocument

# Collapse to the very end of the document and append a new paragraph
# containing "3", matching the centered / size-100 (50pt) formatting
# used by the existing "1" / "2" paragraphs.
$endRange = $d.Content
$endRange.Collapse(0)  # wdCollapseEnd
$endRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newPara.Alignment = 1  # wdAlignParagraphCenter
$newPara.Range.Text = "3"
$newPara.Range.Font.Size = 50
